$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 171, pushing existing rows 171-214 down to 172-215.
$ws.Rows.Item(171).Insert()

# Populate the new row 171 with the same constant columns as its neighbours,
# and the new data point for this market/date.
$ws.Cells.Item(171, 1).Value = 5
$ws.Cells.Item(171, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(171, 3).Value = "Maule"
$ws.Cells.Item(171, 4).Value = 45015
$ws.Cells.Item(171, 5).Value = 7
$ws.Cells.Item(171, 6).Value = 100112030
$ws.Cells.Item(171, 7).Value = "Poroto granado"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Primera"
$ws.Cells.Item(171, 10).Value = 200
$ws.Cells.Item(171, 11).Value = 30000
$ws.Cells.Item(171, 12).Value = 30000
$ws.Cells.Item(171, 13).Value = 30000
$ws.Cells.Item(171, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(171, 15).Value = "Región del Maule"
$ws.Cells.Item(171, 16).Value = 1200
$ws.Cells.Item(171, 17).Value = 25
$ws.Cells.Item(171, 18).Value = "Hortaliza"
